# Update the "Förändrad" (Changed) date column C for rows 2-201
# from serial date 45186 (2023-09-17) to 45188 (2023-09-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 201; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45186) {
        $cell.Value = 45188
    }
}
